$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-sequence the "Periodo Mora" column (E16:E22) into chronological order
# (2307, 2308, 2309, 2310, 2311, 2312, 2401) instead of the previous
# (2401, 2312, 2311, 2310, 2309, 2308, 2307), and swap the "Valor Mora"
# amounts in F16/F22 so the value stays tied to the correct period.

$periods = @("2307", "2308", "2309", "2310", "2311", "2312", "2401")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}

$ws.Range("F16").Value = 46400
$ws.Range("F22").Value = 40000
